# Update "想去人数" (F column) counters across the four sheets to reflect
# the latest scrape snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 345
$ws1.Range("F7").Value  = 3934
$ws1.Range("F9").Value  = 790
$ws1.Range("F10").Value = 2376
$ws1.Range("F11").Value = 370
$ws1.Range("F12").Value = 53
$ws1.Range("F13").Value = 237
$ws1.Range("F14").Value = 762
$ws1.Range("F15").Value = 215
$ws1.Range("F17").Value = 3147
$ws1.Range("F18").Value = 325
$ws1.Range("F21").Value = 354
$ws1.Range("F22").Value = 248
$ws1.Range("F23").Value = 54
$ws1.Range("F24").Value = 283

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value  = 106
$ws2.Range("F10").Value = 101
$ws2.Range("F17").Value = 46
$ws2.Range("F22").Value = 81

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 835
$ws3.Range("F4").Value = 2129
$ws3.Range("F5").Value = 351
$ws3.Range("F6").Value = 21

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 835
$ws4.Range("F4").Value  = 2129
$ws4.Range("F5").Value  = 351
$ws4.Range("F10").Value = 345
$ws4.Range("F16").Value = 21
$ws4.Range("F18").Value = 3934
$ws4.Range("F21").Value = 106
$ws4.Range("F23").Value = 101
$ws4.Range("F24").Value = 790
$ws4.Range("F25").Value = 2376
$ws4.Range("F26").Value = 370
$ws4.Range("F27").Value = 53
$ws4.Range("F29").Value = 237
$ws4.Range("F30").Value = 762
$ws4.Range("F31").Value = 215
$ws4.Range("F35").Value = 325
$ws4.Range("F40").Value = 354
$ws4.Range("F41").Value = 248
$ws4.Range("F42").Value = 54
$ws4.Range("F44").Value = 46
$ws4.Range("F49").Value = 81
$ws4.Range("F50").Value = 283
